$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 (04Okami_Footscray and 05Okami_Brunswick), shifting rows up.
$ws.Range("A4:A5").EntireRow.Delete()

# Update view state: scroll position and active cell selection.
$ws.Range("A4").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
